# "change tracing strategy and save wallet labels"
#
# The token history sheet (Date / USDValue) gets one more observation:
# a new row 3 holding the date "2024-10-05" in column A (column B is left
# blank for this row, matching the existing row 2 pattern).
#
# Note: a plain `$ws.Range("A3").Value = "2024-10-05"` would let Excel's
# input parser auto-recognize the date-shaped string and store it as a
# numeric date serial (with a date number format) instead of literal text.
# The source workbook stores these dates as plain text, so we force the
# cell into Text format before assigning the value, then restore the
# cell's style to Normal (the new cell keeps no explicit formatting, same
# as its sibling A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A3")
$cell.NumberFormat = "@"          # treat the input as text, not a date
$cell.Value = "2024-10-05"        # new wallet-trace date label
$cell.Style = "Normal"            # drop the now-unneeded explicit format
